$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 18, shifting existing rows 18-21 down to 19-22.
$ws.Rows(18).Insert()

# Populate the newly inserted row 18 with this week's data (same
# market/category template as the surrounding rows, with updated date
# and price figures).
$ws.Cells.Item(18, 1).Value = 8
$ws.Cells.Item(18, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 44798
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = 100112026
$ws.Cells.Item(18, 7).Value = "Haba"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 400
$ws.Cells.Item(18, 11).Value = 10500
$ws.Cells.Item(18, 12).Value = 11000
$ws.Cells.Item(18, 13).Value = 10750
$ws.Cells.Item(18, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(18, 16).Value = 430
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
